$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slides")

# Helper: column G holds date-LOOKING text ("YYYY-MM-DD") that must stay plain
# text (matching the rest of the sheet), not get auto-converted to a real date
# serial by Excels input parser. Forcing the cell to Text format first keeps
# the literal string while leaving borders/alignment/etc. untouched.

# --- Update row 2 (S12 / SCOT-HEART: CAC e decisao) ---
$ws.Cells.Item(2,5).Value = "Atualizado"
$ws.Cells.Item(2,7).NumberFormat = "@"
$ws.Cells.Item(2,7).Value = "2026-01-25"
$ws.Cells.Item(2,8).Value = "Figura incluída (KM esquemático) + citação SCOT-HEART 10y corrigida; padding ajustado."
$ws.Cells.Item(2,9).Value = "Lancet 2025 (SCOT-HEART 10y; 10.1016/S0140-6736(24)01899-5); Circulation 2020 (LAP)"
$ws.Cells.Item(2,10).Value = "-"

# --- Update row 18 (S25 / GRADE discriminacao do PREVENT) ---
$ws.Cells.Item(18,5).Value = "Atualizado"
$ws.Cells.Item(18,6).Value = "P2"
$ws.Cells.Item(18,7).NumberFormat = "@"
$ws.Cells.Item(18,7).Value = "2026-01-25"
$ws.Cells.Item(18,8).Value = "Título sem destaque em dourado; tipografia alinhada ao padrão."
$ws.Cells.Item(18,9).Value = "PREVENT/validações externas (referências no slide)"
$ws.Cells.Item(18,10).Value = "-"

# --- Append new rows 20-29 (S50-S59, bloco Metas / Core GRADE / Apendice / Fechamento) ---
# Row 20
$ws.Cells.Item(20,1).Value = "GRADE"
$ws.Cells.Item(20,2).Value = "S50"
$ws.Cells.Item(20,3).Value = "Metas lipídicas (abertura)"
$ws.Cells.Item(20,4).Value = "Metas"
$ws.Cells.Item(20,5).Value = "Novo"
$ws.Cells.Item(20,6).Value = "P0"
$ws.Cells.Item(20,7).NumberFormat = "@"
$ws.Cells.Item(20,7).Value = "2026-01-25"
$ws.Cells.Item(20,8).Value = "Slide de abertura do bloco de metas (SBC 2025)"
$ws.Cells.Item(20,9).Value = "SBC 2025 (10.36660/abc.20250640)"
$ws.Cells.Item(20,10).Value = "-"

# Row 21
$ws.Cells.Item(21,1).Value = "GRADE"
$ws.Cells.Item(21,2).Value = "S51"
$ws.Cells.Item(21,3).Value = "Metas por categoria de risco (LDL/não-HDL/ApoB)"
$ws.Cells.Item(21,4).Value = "Metas"
$ws.Cells.Item(21,5).Value = "Novo"
$ws.Cells.Item(21,6).Value = "P0"
$ws.Cells.Item(21,7).NumberFormat = "@"
$ws.Cells.Item(21,7).Value = "2026-01-25"
$ws.Cells.Item(21,8).Value = "Tabela-resumo com alvos e notas práticas"
$ws.Cells.Item(21,9).Value = "SBC 2025 (Tabela 3.1; 10.36660/abc.20250640)"
$ws.Cells.Item(21,10).Value = "-"

# Row 22
$ws.Cells.Item(22,1).Value = "GRADE"
$ws.Cells.Item(22,2).Value = "S52"
$ws.Cells.Item(22,3).Value = "Força e certeza das metas (conexão com GRADE)"
$ws.Cells.Item(22,4).Value = "Metas"
$ws.Cells.Item(22,5).Value = "Novo"
$ws.Cells.Item(22,6).Value = "P0"
$ws.Cells.Item(22,7).NumberFormat = "@"
$ws.Cells.Item(22,7).Value = "2026-01-25"
$ws.Cells.Item(22,8).Value = "Tabela força/certeza + explicação por domínios/EtD"
$ws.Cells.Item(22,9).Value = "SBC 2025 (GRADE)"
$ws.Cells.Item(22,10).Value = "-"

# Row 23
$ws.Cells.Item(23,1).Value = "GRADE"
$ws.Cells.Item(23,2).Value = "S53"
$ws.Cells.Item(23,3).Value = "Base de evidência: CTT + trials + atualização (VESALIUS-CV)"
$ws.Cells.Item(23,4).Value = "Metas"
$ws.Cells.Item(23,5).Value = "Novo"
$ws.Cells.Item(23,6).Value = "P0"
$ws.Cells.Item(23,7).NumberFormat = "@"
$ws.Cells.Item(23,7).Value = "2026-01-25"
$ws.Cells.Item(23,8).Value = "Resumo fundação + trials + 'living evidence'"
$ws.Cells.Item(23,9).Value = "SBC 2025 refs; NEJMoa2514428"
$ws.Cells.Item(23,10).Value = "-"

# Row 24
$ws.Cells.Item(24,1).Value = "GRADE"
$ws.Cells.Item(24,2).Value = "S54"
$ws.Cells.Item(24,3).Value = "Inconsistência: quando rebaixar?"
$ws.Cells.Item(24,4).Value = "Core GRADE"
$ws.Cells.Item(24,5).Value = "Novo"
$ws.Cells.Item(24,6).Value = "P0"
$ws.Cells.Item(24,7).NumberFormat = "@"
$ws.Cells.Item(24,7).Value = "2026-01-25"
$ws.Cells.Item(24,8).Value = "Slide didático de inconsistência (outlier/explicação)"
$ws.Cells.Item(24,9).Value = "BMJ Core GRADE 3 (10.1136/bmj-2024-081905)"
$ws.Cells.Item(24,10).Value = "-"

# Row 25
$ws.Cells.Item(25,1).Value = "GRADE"
$ws.Cells.Item(25,2).Value = "S55"
$ws.Cells.Item(25,3).Value = "Viés de publicação: como suspeitar?"
$ws.Cells.Item(25,4).Value = "Core GRADE"
$ws.Cells.Item(25,5).Value = "Novo"
$ws.Cells.Item(25,6).Value = "P0"
$ws.Cells.Item(25,7).NumberFormat = "@"
$ws.Cells.Item(25,7).Value = "2026-01-25"
$ws.Cells.Item(25,8).Value = "Sinais e ações práticas no GRADE"
$ws.Cells.Item(25,9).Value = "BMJ Core GRADE 4 (10.1136/bmj-2024-083864)"
$ws.Cells.Item(25,10).Value = "-"

# Row 26
$ws.Cells.Item(26,1).Value = "GRADE"
$ws.Cells.Item(26,2).Value = "S56"
$ws.Cells.Item(26,3).Value = "Divergências entre diretrizes (alto nível)"
$ws.Cells.Item(26,4).Value = "Metas"
$ws.Cells.Item(26,5).Value = "Novo"
$ws.Cells.Item(26,6).Value = "P1"
$ws.Cells.Item(26,7).NumberFormat = "@"
$ws.Cells.Item(26,7).Value = "2026-01-25"
$ws.Cells.Item(26,8).Value = "Tabela comparativa SBC vs ESC/EAS vs ACC vs AACE"
$ws.Cells.Item(26,9).Value = "SBC 2025; ESC/EAS 2019; ACC ECDP 2022; AACE visual guide"
$ws.Cells.Item(26,10).Value = "-"

# Row 27
$ws.Cells.Item(27,1).Value = "GRADE"
$ws.Cells.Item(27,2).Value = "S57"
$ws.Cells.Item(27,3).Value = "Take-home: meta é decisão (EtD)"
$ws.Cells.Item(27,4).Value = "Metas"
$ws.Cells.Item(27,5).Value = "Novo"
$ws.Cells.Item(27,6).Value = "P0"
$ws.Cells.Item(27,7).NumberFormat = "@"
$ws.Cells.Item(27,7).Value = "2026-01-25"
$ws.Cells.Item(27,8).Value = "Fechamento do bloco + ponte para apêndice"
$ws.Cells.Item(27,9).Value = "SBC 2025 + GRADE/EtD"
$ws.Cells.Item(27,10).Value = "-"

# Row 28
$ws.Cells.Item(28,1).Value = "GRADE"
$ws.Cells.Item(28,2).Value = "S58"
$ws.Cells.Item(28,3).Value = "Apêndice (divisor)"
$ws.Cells.Item(28,4).Value = "Apêndice"
$ws.Cells.Item(28,5).Value = "Novo"
$ws.Cells.Item(28,6).Value = "P2"
$ws.Cells.Item(28,7).NumberFormat = "@"
$ws.Cells.Item(28,7).Value = "2026-01-25"
$ws.Cells.Item(28,8).Value = "Slide divisor para manter PREVENT/CAC no fim"
$ws.Cells.Item(28,9).Value = "-"
$ws.Cells.Item(28,10).Value = "-"

# Row 29
$ws.Cells.Item(29,1).Value = "GRADE"
$ws.Cells.Item(29,2).Value = "S59"
$ws.Cells.Item(29,3).Value = "Encerramento com poesia (Camões)"
$ws.Cells.Item(29,4).Value = "Fechamento"
$ws.Cells.Item(29,5).Value = "Novo"
$ws.Cells.Item(29,6).Value = "P2"
$ws.Cells.Item(29,7).NumberFormat = "@"
$ws.Cells.Item(29,7).Value = "2026-01-25"
$ws.Cells.Item(29,8).Value = "Slide final (domínio público)"
$ws.Cells.Item(29,9).Value = "Camões (domínio público)"
$ws.Cells.Item(29,10).Value = "-"

